$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.351.50"
$ws.Range("D3").Value = "1.570.17"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.43"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.34"
$ws.Range("E8").Value = "  -3.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.74"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0895"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "1.795.91"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "1.564.93"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "28.358.31"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.36"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.75"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.10"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  +3.53%  "
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.09"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("D35").Value = "1.379.11"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.521"
$ws.Range("E41").Value = "  -2.45%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.91"
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.784"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0470"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("E46").Value = "  -3.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.28"
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.918"
$ws.Range("E48").Value = "  -6.26%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.708.60"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "85.32"
$ws.Range("E51").Value = "  -0.66%  "
